# Bug fix: interview, wrk exp & minor error msg
#
# Row 2 previously listed "clinical ophthalmic practice msc" with its URL.
# Replace it with the "audiology" program (Advanced Audiology MSc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "https://www.ucl.ac.uk/prospective-students/graduate/taught-degrees/advanced-audiology-msc"
$ws.Range("A2").Value = "audiology"

# Restore the cursor/selection position as recorded in the saved file.
$ws.Range("I13").Select()
